$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-17
# from serial date 45186 to serial date 45188, preserving existing formatting.
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
